# Updating fly and pupa data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Part 1: swap the males/females (D/E) columns for a batch of existing
# rows in the vial=341 block. These rows had D and E entered in the wrong
# order; fix by swapping the two values. ---
$swapRows = @(149,151,153,154,156,157,158,160,161,162,164,165,166,167,168,169,170,172,173,175,176,178,183,184,185,187,189,190,191,193,194,198,199)

foreach ($r in $swapRows) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $dVal
}

# --- Part 2: append newly-collected vial=365 observations as rows 205-233 ---
$newRows = @(
    @(1,  "conditioned",   365, 1, 2),
    @(1,  "unconditioned", 365, 0, 0),
    @(2,  "conditioned",   365, 1, 1),
    @(2,  "unconditioned", 365, 0, 0),
    @(3,  "conditioned",   365, 0, 1),
    @(3,  "unconditioned", 365, 0, 0),
    @(4,  "conditioned",   365, 1, 4),
    @(4,  "unconditioned", 365, 6, 4),
    @(5,  "conditioned",   365, 1, 0),
    @(5,  "unconditioned", 365, 3, 3),
    @(6,  "conditioned",   365, 2, 3),
    @(6,  "unconditioned", 365, 0, 1),
    @(7,  "conditioned",   365, 0, 0),
    @(7,  "unconditioned", 365, 3, 4),
    @(8,  "conditioned",   365, 4, 1),
    @(8,  "unconditioned", 365, 0, 1),
    @(9,  "conditioned",   365, 4, 2),
    @(9,  "unconditioned", 365, 1, 0),
    @(10, "conditioned",   365, 4, 2),
    @(10, "unconditioned", 365, 4, 2),
    @(11, "conditioned",   365, 2, 1),
    @(11, "unconditioned", 365, 0, 0),
    @(12, "conditioned",   365, 0, 1),
    @(12, "unconditioned", 365, 1, 5),
    @(13, "conditioned",   365, 0, 0),
    @(13, "unconditioned", 365, 0, 0),
    @(14, "conditioned",   365, 3, 4),
    @(14, "unconditioned", 365, 0, 0),
    @(15, "unconditioned", 365, 0, 1)
)

$startRow = 205
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# --- Part 3: update the view so the new bottom of the sheet is in frame ---
$ws.Range("E233").Select()
$excel.ActiveWindow.Zoom = 161
$excel.ActiveWindow.ScrollRow = 218
